# Apply updated Price (D) and Volume(1h) (E) values for the crypto symbol list.
# Cells are plain text (inlineStr) in the source workbook, so we force a text
# number format while assigning, then restore the default "Normal" style so no
# extraneous cell style is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "328.85"
Set-TextValue $ws.Range("E2") "-0.59%"

Set-TextValue $ws.Range("D3") "43.23"
Set-TextValue $ws.Range("E3") "2.75%"

Set-TextValue $ws.Range("D4") "5.620"
Set-TextValue $ws.Range("E4") "-0.48%"

Set-TextValue $ws.Range("D5") "0.08205"
Set-TextValue $ws.Range("E5") "-1.60%"

Set-TextValue $ws.Range("D6") "8.763"
Set-TextValue $ws.Range("E6") "-0.41%"

Set-TextValue $ws.Range("D7") "4.456"
Set-TextValue $ws.Range("E7") "-1.43%"

Set-TextValue $ws.Range("D8") "1.909"
Set-TextValue $ws.Range("E8") "-6.27%"

Set-TextValue $ws.Range("D9") "2.853"
Set-TextValue $ws.Range("E9") "-3.87%"

Set-TextValue $ws.Range("D10") "0.9460"
Set-TextValue $ws.Range("E10") "1.79%"

Set-TextValue $ws.Range("E11") "-5.60%"

Set-TextValue $ws.Range("E12") "-2.49%"

Set-TextValue $ws.Range("D13") "0.09756"
Set-TextValue $ws.Range("E13") "3.76%"

Set-TextValue $ws.Range("D14") "0.04443"
Set-TextValue $ws.Range("E14") "13.28%"

Set-TextValue $ws.Range("E15") "0.80%"

Set-TextValue $ws.Range("D16") "0.001286"
Set-TextValue $ws.Range("E16") "-0.62%"

Set-TextValue $ws.Range("D17") "0.006095"
Set-TextValue $ws.Range("E17") "-1.85%"

Set-TextValue $ws.Range("E18") "1.11%"

Set-TextValue $ws.Range("D20") "8.792"
Set-TextValue $ws.Range("E20") "5.49%"

Set-TextValue $ws.Range("D21") "0.1372"
Set-TextValue $ws.Range("E21") "0.92%"

Set-TextValue $ws.Range("D22") "0.2730"
Set-TextValue $ws.Range("E22") "11.34%"

Set-TextValue $ws.Range("D23") "0.04406"
Set-TextValue $ws.Range("E23") "-0.17%"

Set-TextValue $ws.Range("D24") "0.001246"
Set-TextValue $ws.Range("E24") "-0.33%"

Set-TextValue $ws.Range("D25") "0.004391"
Set-TextValue $ws.Range("E25") "-0.01%"

Set-TextValue $ws.Range("D26") "0.0001237"
Set-TextValue $ws.Range("E26") "3.24%"

Set-TextValue $ws.Range("D27") "0.0004013"
Set-TextValue $ws.Range("E27") "31.78%"

Set-TextValue $ws.Range("D39") "0.02780"
Set-TextValue $ws.Range("E39") "-1.77%"

Set-TextValue $ws.Range("D40") "0.05722"
Set-TextValue $ws.Range("E40") "3.03%"

Set-TextValue $ws.Range("D41") "0.007922"
Set-TextValue $ws.Range("E41") "1.74%"

Set-TextValue $ws.Range("D42") "0.009929"
Set-TextValue $ws.Range("E42") "11.32%"

Set-TextValue $ws.Range("D43") "0.1417"
Set-TextValue $ws.Range("E43") "-1.89%"

Set-TextValue $ws.Range("D44") "0.002104"
Set-TextValue $ws.Range("E44") "-6.14%"

Set-TextValue $ws.Range("D45") "0.009738"
Set-TextValue $ws.Range("E45") "-17.04%"

Set-TextValue $ws.Range("D46") "0.00007287"
Set-TextValue $ws.Range("E46") "4.09%"

Set-TextValue $ws.Range("E47") "0.87%"

Set-TextValue $ws.Range("D48") "0.003393"
Set-TextValue $ws.Range("E48") "6.99%"

Set-TextValue $ws.Range("D49") "0.002284"
Set-TextValue $ws.Range("E49") "0.15%"

Set-TextValue $ws.Range("D50") "0.00002113"
Set-TextValue $ws.Range("E50") "0.87%"

Set-TextValue $ws.Range("D51") "0.0002012"
Set-TextValue $ws.Range("E51") "0.87%"
